# Update countries & provincias Spain
# - Refresh the "Datos actualizados" timestamp string
# - Refresh COVID figures for Belgica (row 31) and Kazajistan (row 43)
# - Bahamas overtakes Congo and Trinidad yTobago in total cases, so the
#   three rows (128-130) are re-sorted (descending by "Casos totales"):
#     row128: Bahamas (new, higher figures)
#     row129: Congo   (previous row-128 figures, unchanged)
#     row130: Trinidad yTobago (previous row-129 figures, unchanged)
# - Refresh Recuperados/Casos activos for Butan (row 187)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Footer timestamp (row 1)
$ws.Range("A1").Value = "Datos actualizados a 13 de Octubre de 2020 a las 05:41"

# Row 31 - Belgica
$ws.Cells.Item(31, 2).Value = 165880
$ws.Cells.Item(31, 3).Value = 3622
$ws.Cells.Item(31, 4).Value = 20324
$ws.Cells.Item(31, 5).Value = 135345
$ws.Cells.Item(31, 7).Value = 20
$ws.Cells.Item(31, 8).Value = 10211

# Row 43 - Kazajistan
$ws.Cells.Item(43, 2).Value = 108901
$ws.Cells.Item(43, 3).Value = 70
$ws.Cells.Item(43, 4).Value = 104203
$ws.Cells.Item(43, 5).Value = 2952

# Row 128 - now Bahamas (moved up, updated figures)
$ws.Cells.Item(128, 1).Value = "Bahamas"
$ws.Cells.Item(128, 2).Value = 5163
$ws.Cells.Item(128, 3).Value = 85
$ws.Cells.Item(128, 4).Value = 2978
$ws.Cells.Item(128, 5).Value = 2077
$ws.Cells.Item(128, 7).Value = 1
$ws.Cells.Item(128, 8).Value = 108

# Row 129 - now Congo (previous row-128 data, unchanged)
$ws.Cells.Item(129, 1).Value = "Congo"
$ws.Cells.Item(129, 2).Value = 5118
$ws.Cells.Item(129, 3).Value = 0
$ws.Cells.Item(129, 4).Value = 3887
$ws.Cells.Item(129, 5).Value = 1141
$ws.Cells.Item(129, 7).Value = 0
$ws.Cells.Item(129, 8).Value = 90

# Row 130 - now Trinidad yTobago (previous row-129 data, unchanged)
$ws.Cells.Item(130, 1).Value = "Trinidad yTobago"
$ws.Cells.Item(130, 2).Value = 5116
$ws.Cells.Item(130, 3).Value = 0
$ws.Cells.Item(130, 4).Value = 3303
$ws.Cells.Item(130, 5).Value = 1721
$ws.Cells.Item(130, 7).Value = 0
$ws.Cells.Item(130, 8).Value = 92

# Row 187 - Butan
$ws.Cells.Item(187, 4).Value = 291
$ws.Cells.Item(187, 5).Value = 18
